# pdf results generation, new add to database script, and frontend work
#
# The sample-questions sheet gains a new leading "Exam Name" column: a new
# column is inserted before column A, with a header ("Exam Name") in row 1
# and the exam's name ("mock exam1") in row 2. All the pre-existing data
# (Section/Question/Answer columns) shifts one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; existing columns A-F become B-G.
$ws.Columns("A").Insert()

# Populate the new column's header and the exam name value.
$ws.Range("A1").Value = "Exam Name"
$ws.Range("A2").Value = "mock exam1"

# Size the new column similarly to the other label column (narrower than the
# wide question/answer columns).
$ws.Columns("A").ColumnWidth = 12.43

# Match the saved selection/active cell state.
$ws.Range("A23").Select() | Out-Null
